$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45190 = 2023-09-21) for
# every data row (2-224). The workbook was refreshed and the value needs to
# move forward two days to 45192 (2023-09-23) for all of them.
for ($r = 2; $r -le 224; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value = 45192
    }
}
